# Applies the "Player Info" sheet addition and MATCH_CARD_LINK -> MATCH_CODE
# transformation described in the commit diff.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1. Insert a brand-new "Player Info" sheet as the very first sheet
# -----------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $playerInfo.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the bold / bordered / center-top-aligned header style used by the
# existing sheets' header rows (style index 1 in the original workbook).
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# "ID" (4574) is purely numeric text in the source data, so prefix it with
# an apostrophe -- same as typing it into Excel -- to keep it stored as
# text instead of silently becoming a number.
$values = @("'4574", "Muthuthanthrige Vishwa Thilina Fernando", "Right Handed", "Left Arm Medium Fast")
for ($i = 0; $i -lt $values.Length; $i++) {
    $playerInfo.Cells.Item(2, $i + 1).Value = $values[$i]
}

# -----------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace the URLs with
#    just the bare match-code numbers on the "ODI Batting" sheet
#    (column D) and the "ODI Bowling" sheet (column B). The codes stay
#    text, so prefix with an apostrophe as above.
# -----------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"
for ($r = 2; $r -le 9; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $text = $cell.Text
    if ($text -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
    }
}

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"
for ($r = 2; $r -le 9; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $text = $cell.Text
    if ($text -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
    }
}

Write-Output "done"
